$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("exp_9")

# Fill in column G (msg_overhead_B) for rows 2-16 with numeric values
$ws.Range("G2").Value = 28
$ws.Range("G3").Value = 28
$ws.Range("G4").Value = 28
$ws.Range("G5").Value = 28
$ws.Range("G6").Value = 28
$ws.Range("G7").Value = 28
$ws.Range("G8").Value = 28
$ws.Range("G9").Value = 28
$ws.Range("G10").Value = 28
$ws.Range("G11").Value = 28
$ws.Range("G12").Value = 28
$ws.Range("G13").Value = 48
$ws.Range("G14").Value = 68
$ws.Range("G15").Value = 128
$ws.Range("G16").Value = 248

# Mark the fragmentation note next to row 13
$ws.Range("H13").Value = "* begin fragmentation"

# Remaining rows (fragmented past the MTU) are marked "na"
$ws.Range("G17").Value = "na"
$ws.Range("G18").Value = "na"
$ws.Range("G19").Value = "na"
$ws.Range("G20").Value = "na"
$ws.Range("G21").Value = "na"
$ws.Range("G22").Value = "na"

# Update the selection shown on this sheet
$ws.Range("H18").Select()

# Make exp_9 the active (selected) tab, replacing exp_8
$ws.Activate()
